$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (42) describing the "glucosidase" target protein, following
# the existing table layout (columns A:N).
$seq = "MRFPSIFTAVLFAASSALAAPVNTTTEDETAQIPAEAVIGYSDLEGDFDVAVLPFSNSTNNGLLFINTTIASIAAKEEGVSLDKREEGEPKASIPSSASVQLDSYNYDGSTFSGKIYVKNIAYSKKVTVIYADGSDNWNNNGNTIAASYSAPISGSNYEYWTFSASINGIKEFYIKYEVSGKTYYDNNNSANYQVSTSKPTTTTATATTTTAPSTSTTTPPSRSEPATFPTGNSTISSWIKKQEGISRFAMLRNINPPGSATGFIAASLSTAGPDYYYAWTRDAALTSNVIVYEYNTTLSGNKTILNVLKDYVTFSVKTQSTSTVCNCLGEPKFNPDASGYTGAWGRPQNDGPAERATTFILFADSYLTQTKDASYVTGTLKPAIFKDLDYVVNVWSNGCFDLWEEVNGVHFYTLMVMRKGLLLGADFAKRNGDSTRASTYSSTASTIANKISSFWVSSNNWIQVSQSVTGGVSKKGLDVSTLLAANLGSVDDGFFTPGSEKILATAVAVEDSFASLYPINKNLPSYLGNSIGRYPEDTYNGNGNSQGNSWFLAVTGYAELYYRAIKEWIGNGGVTVSSISLPFFKKFDSSATSGKKYTVGTSDFNNLAQNIALAADRFLSTVQLHAHNNGSLAEEFDRTTGLSTGARDLTWSHASLITASYAKAGAPAA"
$sp = "MRFPSIFTAVLFAASSALAAPVNTTTEDETAQIPAEAVIGYSDLEGDFDVAVLPFSNSTNNGLLFINTTIASIAAKEEGVSLDKREEGEPK"

$row = 42
$ws.Cells.Item($row, 1).Value = "glucosidase"
$ws.Cells.Item($row, 2).Value = "glucosidase"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 5
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = "e"
$ws.Cells.Item($row, 11).Value = $seq
$ws.Cells.Item($row, 12).Value = 670
$ws.Cells.Item($row, 13).Value = $sp
$ws.Cells.Item($row, 14).Value = 91

# Column M ("sp sequence") is highlighted the same way as the other rows in
# that column - orange fill (matches existing style used throughout M2:M41).
$ws.Cells.Item($row, 13).Interior.Color = 49407

# Reflect the author's last selection before saving.
$ws.Range("M28").Select() | Out-Null
